# Scheduled market-data refresh: pushes newly-fetched Universalis price
# snapshots (currentAveragePrice / NQ / HQ) and recomputed Leve profit
# figures into each crafting-job worksheet. Mirrors the output of the
# repo's scheduled runner job ("chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1080.3043
$ws.Range("I11").Value = 1080.3043
$ws.Range("K11").Value = 1080.3043
$ws.Range("M11").Value = -940.3043
$ws.Range("H40").Value = 2705.5789
$ws.Range("I40").Value = 4560
$ws.Range("J40").Value = 2043.2858
$ws.Range("K40").Value = 4560
$ws.Range("L40").Value = 2043.2858
$ws.Range("M40").Value = -4385
$ws.Range("N40").Value = -2393.2858
$ws.Range("I64").Value = 3404.0908
$ws.Range("J64").Value = 3900
$ws.Range("K64").Value = 3404.0908
$ws.Range("L64").Value = 3900
$ws.Range("M64").Value = -3156.0908
$ws.Range("N64").Value = -4396
$ws.Range("I67").Value = 3404.0908
$ws.Range("J67").Value = 3900
$ws.Range("K67").Value = 3404.0908
$ws.Range("L67").Value = 3900
$ws.Range("M67").Value = -2546.0908
$ws.Range("N67").Value = -5616
$ws.Range("H74").Value = 3994.0588
$ws.Range("I74").Value = 3785.7144
$ws.Range("J74").Value = 4139.9
$ws.Range("K74").Value = 3785.7144
$ws.Range("L74").Value = 4139.9
$ws.Range("M74").Value = -2849.7144
$ws.Range("N74").Value = -6011.9
$ws.Range("H76").Value = 3729.7026
$ws.Range("I76").Value = 3612.9666
$ws.Range("K76").Value = 3612.9666
$ws.Range("M76").Value = -3297.9666
$ws.Range("H77").Value = 3994.0588
$ws.Range("I77").Value = 3785.7144
$ws.Range("J77").Value = 4139.9
$ws.Range("K77").Value = 18928.572
$ws.Range("L77").Value = 20699.5
$ws.Range("M77").Value = -14248.572
$ws.Range("N77").Value = -30059.5
$ws.Range("H79").Value = 3729.7026
$ws.Range("I79").Value = 3612.9666
$ws.Range("K79").Value = 3612.9666
$ws.Range("M79").Value = -2520.9666
$ws.Range("H100").Value = 2408
$ws.Range("I100").Value = 416.66666
$ws.Range("K100").Value = 416.66666
$ws.Range("M100").Value = 124.33334
$ws.Range("H113").Value = 2914.8572
$ws.Range("I113").Value = 2750.8
$ws.Range("J113").Value = 3006
$ws.Range("K113").Value = 2750.8
$ws.Range("L113").Value = 3006
$ws.Range("M113").Value = 503.1999999999998
$ws.Range("N113").Value = -9514
$ws.Range("H116").Value = 4362.727
$ws.Range("I116").Value = 4501.25
$ws.Range("J116").Value = 3993.3333
$ws.Range("K116").Value = 4501.25
$ws.Range("L116").Value = 3993.3333
$ws.Range("M116").Value = -1059.25
$ws.Range("N116").Value = -10877.3333
$ws.Range("H132").Value = 1761.6078
$ws.Range("I132").Value = 1557.1951
$ws.Range("J132").Value = 2599.7
$ws.Range("K132").Value = 4671.5853
$ws.Range("L132").Value = 7799.099999999999
$ws.Range("M132").Value = -2141.5853
$ws.Range("N132").Value = -12859.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6121.6
$ws.Range("I61").Value = 4465.5
$ws.Range("K61").Value = 4465.5
$ws.Range("M61").Value = -4253.5
$ws.Range("H81").Value = 35640.25
$ws.Range("J81").Value = 35640.25
$ws.Range("L81").Value = 35640.25
$ws.Range("N81").Value = -37636.25
$ws.Range("H84").Value = 35640.25
$ws.Range("J84").Value = 35640.25
$ws.Range("L84").Value = 106920.75
$ws.Range("N84").Value = -116904.75
$ws.Range("H132").Value = 5918.4595
$ws.Range("I132").Value = 1853.5
$ws.Range("K132").Value = 5560.5
$ws.Range("M132").Value = -3030.5
$ws.Range("H136").Value = 6121.6
$ws.Range("I136").Value = 4465.5
$ws.Range("K136").Value = 13396.5
$ws.Range("M136").Value = -10846.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5051.5884
$ws.Range("I134").Value = 4313.393
$ws.Range("J134").Value = 8496.5
$ws.Range("K134").Value = 12940.179
$ws.Range("L134").Value = 25489.5
$ws.Range("M134").Value = -10405.179
$ws.Range("N134").Value = -30559.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2528861.8
$ws.Range("I58").Value = 4331848
$ws.Range("J58").Value = 4680.933
$ws.Range("K58").Value = 4331848
$ws.Range("L58").Value = 4680.933
$ws.Range("M58").Value = -4331645
$ws.Range("N58").Value = -5086.933
$ws.Range("H99").Value = 2735.6667
$ws.Range("I99").Value = 2827.625
$ws.Range("K99").Value = 2827.625
$ws.Range("M99").Value = -1329.625
$ws.Range("H122").Value = 10516.05
$ws.Range("I122").Value = 3964.7856
$ws.Range("J122").Value = 25802.334
$ws.Range("K122").Value = 11894.3568
$ws.Range("L122").Value = 77407.00199999999
$ws.Range("M122").Value = -9444.356800000001
$ws.Range("N122").Value = -82307.00199999999
$ws.Range("H126").Value = 2735.6667
$ws.Range("I126").Value = 2827.625
$ws.Range("K126").Value = 8482.875
$ws.Range("M126").Value = -6012.875
$ws.Range("H132").Value = 2644.4
$ws.Range("I132").Value = 2438.2778
$ws.Range("K132").Value = 7314.8334
$ws.Range("M132").Value = -4784.8334
$ws.Range("H134").Value = 22497.51
$ws.Range("I134").Value = 38321.605
$ws.Range("K134").Value = 114964.815
$ws.Range("M134").Value = -112429.815
$ws.Range("H136").Value = 2528861.8
$ws.Range("I136").Value = 4331848
$ws.Range("J136").Value = 4680.933
$ws.Range("K136").Value = 12995544
$ws.Range("L136").Value = 14042.799
$ws.Range("M136").Value = -12992994
$ws.Range("N136").Value = -19142.799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 27591.139
$ws.Range("J131").Value = 40565.668
$ws.Range("L131").Value = 121697.004
$ws.Range("N131").Value = -131777.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3506.7742
$ws.Range("I102").Value = 3810.4736
$ws.Range("J102").Value = 3025.9167
$ws.Range("K102").Value = 3810.4736
$ws.Range("L102").Value = 3025.9167
$ws.Range("M102").Value = -2188.4736
$ws.Range("N102").Value = -6269.9167
$ws.Range("H126").Value = 3377.1765
$ws.Range("I126").Value = 1990.2222
$ws.Range("J126").Value = 4937.5
$ws.Range("K126").Value = 5970.6666
$ws.Range("L126").Value = 14812.5
$ws.Range("M126").Value = -3500.6666
$ws.Range("N126").Value = -19752.5
$ws.Range("H132").Value = 5017.533
$ws.Range("I132").Value = 9809.25
$ws.Range("K132").Value = 29427.75
$ws.Range("M132").Value = -26897.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30920
$ws.Range("H122").Value = 4815.075
$ws.Range("I122").Value = 4117.517
$ws.Range("J122").Value = 6654.091
$ws.Range("K122").Value = 12352.551
$ws.Range("L122").Value = 19962.273
$ws.Range("M122").Value = -9902.550999999999
$ws.Range("N122").Value = -24862.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 12000
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12470
$ws.Range("H35").Value = 12000
$ws.Range("J35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("N35").Value = -12580
$ws.Range("H75").Value = 24423.75
$ws.Range("J75").Value = 27398.334
$ws.Range("L75").Value = 27398.334
$ws.Range("N75").Value = -29270.334
$ws.Range("H78").Value = 24423.75
$ws.Range("J78").Value = 27398.334
$ws.Range("L78").Value = 82195.00199999999
$ws.Range("N78").Value = -91555.00199999999
